$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.201.89"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.852.59"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.0000"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +1.36%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "237.77"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  +0.02%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07854"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.01%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3018"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.76"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.42%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08112"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.839.44"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.17%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.174"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("E14").Value = "  -2.33%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "89.50"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "29.223.45"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("E17").Value = "  +1.06%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000007821"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.21"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "235.53"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.9998"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.102.49"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("E23").Value = "  +0.03%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.510"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "162.63"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.866"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1423"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.01"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("E29").Value = "  -1.82%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.406"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.478"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.325"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -4.22%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.009"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05154"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.164"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.66%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7103"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.9968"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("E40").Value = "  +0.94%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.150.95"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +5.26%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9253"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("E43").Value = "  -0.05%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4233"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "70.09"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("E46").Value = "  +0.04%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "103.00"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5296"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.79%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.736"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -3.60%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "9.136"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "6.951"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.77%  "
